$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# D-column values are forced to Text via a leading apostrophe (quote-prefix)
# so numeric-looking strings (e.g. "486.68") are not auto-converted to
# Number cells by Excel; Style is reset to "Normal" afterwards so the
# quote-prefix formatting flag does not linger as a new cell style.

$ws.Range("D2").Value = "`'68.330.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "`'3.928.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("D4").Value = "`'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "`'486.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.44%  "
$ws.Range("D6").Value = "`'148.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.64%  "
$ws.Range("D7").Value = "`'0.631"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.43%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").Value = "  +3.91%  "
$ws.Range("D11").Value = "`'0.0000353"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.33%  "
$ws.Range("D12").Value = "`'43.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "`'10.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.32%  "
$ws.Range("D14").Value = "`'4.549.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").Value = "`'3.929.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "`'20.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").Value = "`'1.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").Value = "`'68.397.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.83%  "
$ws.Range("D21").Value = "`'441.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.61%  "
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("E23").Value = "  +2.26%  "
$ws.Range("D24").Value = "`'88.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("D25").Value = "`'11.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +14.00%  "
$ws.Range("D26").Value = "`'11.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +15.51%  "
$ws.Range("D27").Value = "`'3.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("D28").Value = "`'38.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("D29").Value = "`'5.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.79%  "
$ws.Range("D30").Value = "`'718.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("D31").Value = "`'13.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("E33").Value = "  +3.08%  "
$ws.Range("D34").Value = "`'0.0₃0914"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +15.62%  "
$ws.Range("D35").Value = "`'42.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.26%  "
$ws.Range("D36").Value = "`'6.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +13.00%  "
$ws.Range("D37").Value = "`'59.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.61%  "
$ws.Range("E38").Value = "  -3.06%  "
$ws.Range("D39").Value = "`'0.397"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +17.58%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "`'2.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +14.36%  "
$ws.Range("D42").Value = "`'0.0483"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.68%  "
$ws.Range("D43").Value = "`'3.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.54%  "
$ws.Range("D44").Value = "`'2.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.94%  "
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "`'0.0₆0356"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +46.69%  "
$ws.Range("D48").Value = "`'3.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").Value = "`'145.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("D51").Value = "`'3.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.75%  "
